$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-coerced to a number by Excel
# (losing the literal text formatting, e.g. trailing zeros) are forced to Text
# format before the write, then restored to the default style so no stray
# number-format style lingers on the cell.

$ws.Range("D2").Value = '26.648.93'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '1.854.58'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '264.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.66%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5268'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3251'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06806'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.96'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7849'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07787'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("D13").Value = '1.850.97'
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.034'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007990'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = '26.674.86'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("E21").Value = '  +2.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.482'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.019'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.176'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.693'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.194'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.120'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08729'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04864'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7206'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.876'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.119'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.260'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01796'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.4880'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9026'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.976'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.002'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.697'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4209'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05892'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.049'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1239'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.18'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.8903'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '60.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.37%  '
